$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"23.94933966666666"
$ws.Range("H2").Value = [double]"71.84801899999999"
$ws.Range("I2").Value = [double]"0.003771463579284448"
$ws.Range("J2").Value = [double]"0.003771463579284448"
$ws.Range("K2").Value = [double]"3"
$ws.Range("L2").Value = [double]"1"
$ws.Range("M2").Value = [double]"1.728078666666667"
$ws.Range("N2").Value = [double]"5.184236"
$ws.Range("O2").Value = [double]"0.01495566191401176"
$ws.Range("P2").Value = [double]"0.01495566191401175"
$ws.Range("Q2").Value = [double]"41.38634295872044"
$ws.Range("R2").Value = [double]"372.477086628484"
$ws.Range("S2").Value = [double]"5.640473421278687E-05"
$ws.Range("T2").Value = [double]"5.640473421278687E-05"
$ws.Range("G3").Value = [double]"23.94933966666666"
$ws.Range("H3").Value = [double]"71.84801899999999"
$ws.Range("I3").Value = [double]"0.003771463579284448"
$ws.Range("J3").Value = [double]"0.003771463579284448"
$ws.Range("O3").Value = [double]"0.6765388985664503"
$ws.Range("P3").Value = [double]"0.6765388985664502"
$ws.Range("Q3").Value = [double]"1872.165273725115"
$ws.Range("R3").Value = [double]"16849.48746352603"
$ws.Range("S3").Value = [double]"0.002551541815912583"
$ws.Range("T3").Value = [double]"0.002551541815912583"
$ws.Range("G4").Value = [double]"23.94933966666666"
$ws.Range("H4").Value = [double]"71.84801899999999"
$ws.Range("I4").Value = [double]"0.003771463579284448"
$ws.Range("J4").Value = [double]"0.003771463579284448"
$ws.Range("M4").Value = [double]"35.284818"
$ws.Range("N4").Value = [double]"105.854454"
$ws.Range("O4").Value = [double]"0.3053725613795956"
$ws.Range("P4").Value = [double]"0.3053725613795956"
$ws.Range("Q4").Value = [double]"845.0480913585139"
$ws.Range("R4").Value = [double]"7605.432822226626"
$ws.Range("S4").Value = [double]"0.001151701493355949"
$ws.Range("T4").Value = [double]"0.001151701493355949"
$ws.Range("G5").Value = [double]"23.94933966666666"
$ws.Range("H5").Value = [double]"71.84801899999999"
$ws.Range("I5").Value = [double]"0.003771463579284448"
$ws.Range("J5").Value = [double]"0.003771463579284448"
$ws.Range("K5").Value = [double]"3"
$ws.Range("L5").Value = [double]"1"
$ws.Range("M5").Value = [double]"0.361994"
$ws.Range("N5").Value = [double]"1.085982"
$ws.Range("O5").Value = [double]"0.003132878139942378"
$ws.Range("P5").Value = [double]"0.003132878139942378"
$ws.Range("Q5").Value = [double]"8.669517263295331"
$ws.Range("R5").Value = [double]"78.025655369658"
$ws.Range("S5").Value = [double]"1.181553580312908E-05"
$ws.Range("T5").Value = [double]"1.181553580312908E-05"
$ws.Range("I6").Value = [double]"0.9508087878751565"
$ws.Range("J6").Value = [double]"0.9508087878751567"
$ws.Range("K6").Value = [double]"3"
$ws.Range("L6").Value = [double]"1"
$ws.Range("M6").Value = [double]"1.728078666666667"
$ws.Range("N6").Value = [double]"5.184236"
$ws.Range("O6").Value = [double]"0.01495566191401176"
$ws.Range("P6").Value = [double]"0.01495566191401175"
$ws.Range("Q6").Value = [double]"10433.7474712224"
$ws.Range("R6").Value = [double]"93903.72724100163"
$ws.Range("S6").Value = [double]"0.01421997477633216"
$ws.Range("T6").Value = [double]"0.01421997477633216"
$ws.Range("I7").Value = [double]"0.9508087878751565"
$ws.Range("J7").Value = [double]"0.9508087878751567"
$ws.Range("O7").Value = [double]"0.6765388985664503"
$ws.Range("P7").Value = [double]"0.6765388985664502"
$ws.Range("S7").Value = [double]"0.6432591300963602"
$ws.Range("T7").Value = [double]"0.64325913009636"
$ws.Range("I8").Value = [double]"0.9508087878751565"
$ws.Range("J8").Value = [double]"0.9508087878751567"
$ws.Range("M8").Value = [double]"35.284818"
$ws.Range("N8").Value = [double]"105.854454"
$ws.Range("O8").Value = [double]"0.3053725613795956"
$ws.Range("P8").Value = [double]"0.3053725613795956"
$ws.Range("Q8").Value = [double]"213041.7368615411"
$ws.Range("R8").Value = [double]"1917375.63175387"
$ws.Range("S8").Value = [double]"0.2903509149356652"
$ws.Range("T8").Value = [double]"0.2903509149356651"
$ws.Range("I9").Value = [double]"0.9508087878751565"
$ws.Range("J9").Value = [double]"0.9508087878751567"
$ws.Range("K9").Value = [double]"3"
$ws.Range("L9").Value = [double]"1"
$ws.Range("M9").Value = [double]"0.361994"
$ws.Range("N9").Value = [double]"1.085982"
$ws.Range("O9").Value = [double]"0.003132878139942378"
$ws.Range("P9").Value = [double]"0.003132878139942378"
$ws.Range("Q9").Value = [double]"2185.637757674042"
$ws.Range("R9").Value = [double]"19670.73981906638"
$ws.Range("S9").Value = [double]"0.002978768066799187"
$ws.Range("T9").Value = [double]"0.002978768066799187"
$ws.Range("G10").Value = [double]"285.3476563333333"
$ws.Range("H10").Value = [double]"856.042969"
$ws.Range("I10").Value = [double]"0.0449356144375536"
$ws.Range("J10").Value = [double]"0.04493561443755361"
$ws.Range("K10").Value = [double]"3"
$ws.Range("L10").Value = [double]"1"
$ws.Range("M10").Value = [double]"1.728078666666667"
$ws.Range("N10").Value = [double]"5.184236"
$ws.Range("O10").Value = [double]"0.01495566191401176"
$ws.Range("P10").Value = [double]"0.01495566191401175"
$ws.Range("Q10").Value = [double]"493.1031974929649"
$ws.Range("R10").Value = [double]"4437.928777436684"
$ws.Range("S10").Value = [double]"0.0006720418574264372"
$ws.Range("T10").Value = [double]"0.0006720418574264372"
$ws.Range("G11").Value = [double]"285.3476563333333"
$ws.Range("H11").Value = [double]"856.042969"
$ws.Range("I11").Value = [double]"0.0449356144375536"
$ws.Range("J11").Value = [double]"0.04493561443755361"
$ws.Range("O11").Value = [double]"0.6765388985664503"
$ws.Range("P11").Value = [double]"0.6765388985664502"
$ws.Range("Q11").Value = [double]"22306.1671245013"
$ws.Range("R11").Value = [double]"200755.5041205117"
$ws.Range("S11").Value = [double]"0.0304006910979892"
$ws.Range("T11").Value = [double]"0.0304006910979892"
$ws.Range("G12").Value = [double]"285.3476563333333"
$ws.Range("H12").Value = [double]"856.042969"
$ws.Range("I12").Value = [double]"0.0449356144375536"
$ws.Range("J12").Value = [double]"0.04493561443755361"
$ws.Range("M12").Value = [double]"35.284818"
$ws.Range("N12").Value = [double]"105.854454"
$ws.Range("O12").Value = [double]"0.3053725613795956"
$ws.Range("P12").Value = [double]"0.3053725613795956"
$ws.Range("Q12").Value = [double]"10068.44012044821"
$ws.Range("R12").Value = [double]"90615.96108403393"
$ws.Range("S12").Value = [double]"0.01372210367796168"
$ws.Range("T12").Value = [double]"0.01372210367796168"
$ws.Range("G13").Value = [double]"285.3476563333333"
$ws.Range("H13").Value = [double]"856.042969"
$ws.Range("I13").Value = [double]"0.0449356144375536"
$ws.Range("J13").Value = [double]"0.04493561443755361"
$ws.Range("K13").Value = [double]"3"
$ws.Range("L13").Value = [double]"1"
$ws.Range("M13").Value = [double]"0.361994"
$ws.Range("N13").Value = [double]"1.085982"
$ws.Range("O13").Value = [double]"0.003132878139942378"
$ws.Range("P13").Value = [double]"0.003132878139942378"
$ws.Range("Q13").Value = [double]"103.2941395067287"
$ws.Range("R13").Value = [double]"929.6472555605579"
$ws.Range("S13").Value = [double]"0.0001407778041762908"
$ws.Range("T13").Value = [double]"0.0001407778041762908"
$ws.Range("G14").Value = [double]"3.074321666666667"
$ws.Range("H14").Value = [double]"9.222965"
$ws.Range("I14").Value = [double]"0.0004841341080053326"
$ws.Range("J14").Value = [double]"0.0004841341080053326"
$ws.Range("K14").Value = [double]"3"
$ws.Range("L14").Value = [double]"1"
$ws.Range("M14").Value = [double]"1.728078666666667"
$ws.Range("N14").Value = [double]"5.184236"
$ws.Range("O14").Value = [double]"0.01495566191401176"
$ws.Range("P14").Value = [double]"0.01495566191401175"
$ws.Range("Q14").Value = [double]"5.312669686637778"
$ws.Range("R14").Value = [double]"47.81402717974"
$ws.Range("S14").Value = [double]"7.240546040369407E-06"
$ws.Range("T14").Value = [double]"7.240546040369406E-06"
$ws.Range("G15").Value = [double]"3.074321666666667"
$ws.Range("H15").Value = [double]"9.222965"
$ws.Range("I15").Value = [double]"0.0004841341080053326"
$ws.Range("J15").Value = [double]"0.0004841341080053326"
$ws.Range("O15").Value = [double]"0.6765388985664503"
$ws.Range("P15").Value = [double]"0.6765388985664502"
$ws.Range("Q15").Value = [double]"240.3255515476656"
$ws.Range("R15").Value = [double]"2162.92996392899"
$ws.Range("S15").Value = [double]"0.0003275355561883787"
$ws.Range("T15").Value = [double]"0.0003275355561883786"
$ws.Range("G16").Value = [double]"3.074321666666667"
$ws.Range("H16").Value = [double]"9.222965"
$ws.Range("I16").Value = [double]"0.0004841341080053326"
$ws.Range("J16").Value = [double]"0.0004841341080053326"
$ws.Range("M16").Value = [double]"35.284818"
$ws.Range("N16").Value = [double]"105.854454"
$ws.Range("O16").Value = [double]"0.3053725613795956"
$ws.Range("P16").Value = [double]"0.3053725613795956"
$ws.Range("Q16").Value = [double]"108.47688048179"
$ws.Range("R16").Value = [double]"976.29192433611"
$ws.Range("S16").Value = [double]"0.0001478412726128142"
$ws.Range("T16").Value = [double]"0.0001478412726128142"
$ws.Range("G17").Value = [double]"3.074321666666667"
$ws.Range("H17").Value = [double]"9.222965"
$ws.Range("I17").Value = [double]"0.0004841341080053326"
$ws.Range("J17").Value = [double]"0.0004841341080053326"
$ws.Range("K17").Value = [double]"3"
$ws.Range("L17").Value = [double]"1"
$ws.Range("M17").Value = [double]"0.361994"
$ws.Range("N17").Value = [double]"1.085982"
$ws.Range("O17").Value = [double]"0.003132878139942378"
$ws.Range("P17").Value = [double]"0.003132878139942378"
$ws.Range("Q17").Value = [double]"1.112885997403333"
$ws.Range("R17").Value = [double]"10.01597397663"
$ws.Range("S17").Value = [double]"1.516733163770409E-06"
$ws.Range("T17").Value = [double]"1.516733163770409E-06"
